# #3473 swapped out two properties
# Update "Portfolio Manager Building ID" (col B) and "Gross Area (SF)" (col L)
# values for a couple of rows, and move the active selection to reflect the
# user's last click (L2:L10, active cell L2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPS Data")

# Gross Area (SF) corrected for row 3
$ws.Range("L3").Value = 227440.2

# Portfolio Manager Building ID corrected for rows 6 and 10
$ws.Range("B6").Value = 22482006
$ws.Range("B10").Value = 22482007

# Reflect the resulting selection/active cell in column L
$ws.Range("L2:L10").Select() | Out-Null
